$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Id" column (A) previously held numeric placeholder ids; replace them
# with the corrected text identifiers.
$ws.Range("A2").Value = "JImport1"
$ws.Range("A3").Value = "Jimport2"
